$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: apply values from original row 11
$ws.Range("D2").Value = [datetime]"2021-11-19"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 2688

# Row 3: apply values from original row 12
$ws.Range("D3").Value = [datetime]"2021-11-19"
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("S3").Value = 2250

# Row 4: apply values from original row 17
$ws.Range("D4").Value = [datetime]"2021-10-19"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("S4").Value = 2188

# Row 5: apply values from original row 4
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2688

# Row 6: apply values from original row 5
$ws.Range("D6").Value = [datetime]"2021-11-23"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("Q6").Value = "$/bandeja 8 kilos"
$ws.Range("S6").Value = 2250
$ws.Range("T6").Value = 8

# Row 7: apply values from original row 2
$ws.Range("D7").Value = [datetime]"2021-10-26"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19556
$ws.Range("S7").Value = 2444

# Row 8: apply values from original row 23
$ws.Range("D8").Value = [datetime]"2021-11-16"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("S8").Value = 2312

# Row 9: apply values from original row 26
$ws.Range("D9").Value = [datetime]"2021-11-26"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 21000
$ws.Range("S9").Value = 2625

# Row 10: apply values from original row 7
$ws.Range("D10").Value = [datetime]"2021-12-03"
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 18500
$ws.Range("S10").Value = 2312

# Row 11: apply values from original row 8
$ws.Range("D11").Value = [datetime]"2021-12-03"
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("S11").Value = 2000

# Row 12: apply values from original row 6
$ws.Range("D12").Value = [datetime]"2020-11-27"
$ws.Range("L12").Value = "Primera"
$ws.Range("N12").Value = 2000
$ws.Range("O12").Value = 2100
$ws.Range("P12").Value = 2050
$ws.Range("Q12").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S12").Value = 2050
$ws.Range("T12").Value = 1

# Row 13: apply values from original row 10
$ws.Range("D13").Value = [datetime]"2022-11-25"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22500
$ws.Range("P13").Value = 22250
$ws.Range("S13").Value = 2781

# Row 14: apply values from original row 20
$ws.Range("D14").Value = [datetime]"2022-11-11"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 22000
$ws.Range("O14").Value = 22500
$ws.Range("P14").Value = 22250
$ws.Range("S14").Value = 2781

# Row 15: apply values from original row 21
$ws.Range("D15").Value = [datetime]"2020-11-24"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 2000
$ws.Range("O15").Value = 2100
$ws.Range("P15").Value = 2050
$ws.Range("Q15").Value = "$/kilo (en caja de 14 kilos)"
$ws.Range("S15").Value = 2050
$ws.Range("T15").Value = 1

# Row 16: apply values from original row 15
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("S16").Value = 2438

# Row 17: apply values from original row 16
$ws.Range("D17").Value = [datetime]"2021-11-30"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 2000

# Row 18: apply values from original row 3
$ws.Range("D18").Value = [datetime]"2021-11-09"
$ws.Range("M18").Value = 200

# Row 19: apply values from original row 14
$ws.Range("D19").Value = [datetime]"2021-10-22"
$ws.Range("L19").Value = "Segunda"
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 19000
$ws.Range("P19").Value = 18500
$ws.Range("S19").Value = 2312

# Row 20: apply values from original row 22
$ws.Range("D20").Value = [datetime]"2022-11-15"

# Row 21: apply values from original row 9
$ws.Range("D21").Value = [datetime]"2022-11-29"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 22500
$ws.Range("P21").Value = 22250
$ws.Range("Q21").Value = "$/bandeja 8 kilos"
$ws.Range("S21").Value = 2781
$ws.Range("T21").Value = 8

# Row 22: apply values from original row 13
$ws.Range("D22").Value = [datetime]"2021-11-05"
$ws.Range("L22").Value = "Segunda"
$ws.Range("N22").Value = 19000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19500
$ws.Range("S22").Value = 2438

# Row 23: apply values from original row 18
$ws.Range("D23").Value = [datetime]"2021-10-29"
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("S23").Value = 2438

# Row 24: apply values from original row 19
$ws.Range("D24").Value = [datetime]"2022-11-30"
$ws.Range("M24").Value = 200

# Row 26: apply values from original row 24
$ws.Range("D26").Value = [datetime]"2022-11-08"
$ws.Range("N26").Value = 22000
$ws.Range("O26").Value = 22500
$ws.Range("P26").Value = 22250
$ws.Range("S26").Value = 2781

